# Applies the "reformat architecture and update the tester nodes" edit:
#  1. Updates a handful of status cells on the "Report" sheet:
#       - Bernoulli (row 9) flips from a FAIL message to
#         "SUCCESS (via decomposition)" and gets a used_provider value.
#       - Several "SUCCESS WITH FALLBACK" rows become
#         "SUCCESS (via decomposition)".
#     (The underlying report generator also reshuffles the internal
#     fill/style table order for the FAIL/UNKNOWN rows that are left
#     untouched, but their rendered color - red for FAIL, grey for
#     UNKNOWN - does not actually change, so there is nothing to do for
#     those cells here.)
#  2. Refreshes the aggregated counts/percentages and the generated-on
#     timestamp on the "Data_PieChart" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Report")

# BGR integer equivalent (as exposed through Interior.Color) of the green
# palette color already used by existing SUCCESS cells.
$colorGreen = 4499968   # RGB 00AA44

# --- 1/2: status text updates ---

# Row 9 - Bernoulli now succeeds via decomposition instead of failing.
$ws.Range("C9").Value = "TensorrtExecutionProvider"
$ws.Range("D9").Value = "SUCCESS (via decomposition)"
$ws.Range("D9").Interior.Color = $colorGreen

# Rows whose "SUCCESS WITH FALLBACK" result became "SUCCESS (via decomposition)".
$decompositionCells = @("D15","D29","D55","D56","D88","D89","D123","D127","D137")
foreach ($addr in $decompositionCells) {
    $ws.Range($addr).Value = "SUCCESS (via decomposition)"
    $ws.Range($addr).Interior.Color = $colorGreen
}

# --- 2/2: refresh the pie-chart summary data ---

$pie = $wb.Worksheets.Item("Data_PieChart")
$pie.Range("B2").Value = 103
$pie.Range("C2").Value = 39.3
$pie.Range("B3").Value = 129
$pie.Range("C3").Value = 49.2
$pie.Range("B7").Value = 21
$pie.Range("C7").Value = 8
$pie.Range("B8").Value = "2025-11-18 14:38:03"
